# 🔄 Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in resultado (G) / profit (H) for rows that were pending and now have a result
$updates = @(
    @{ Row = 108; G = "Acierto"; H = 1.1 },
    @{ Row = 109; G = "Fallo";   H = -1 },
    @{ Row = 110; G = "Fallo";   H = -1 },
    @{ Row = 116; G = "Fallo";   H = -1 },
    @{ Row = 117; G = "Fallo";   H = -1 },
    @{ Row = 118; G = "Fallo";   H = -1 },
    @{ Row = 119; G = "Fallo";   H = -1 },
    @{ Row = 124; G = "Acierto"; H = 2 },
    @{ Row = 125; G = "Fallo";   H = -1 },
    @{ Row = 126; G = "Fallo";   H = -1 },
    @{ Row = 130; G = "Fallo";   H = -1 },
    @{ Row = 131; G = "Fallo";   H = -1 }
)

foreach ($u in $updates) {
    $ws.Range("G" + $u.Row).Value = $u.G
    $ws.Range("H" + $u.Row).Value = $u.H
}

# Append new matches to the tracker
$newRows = @(
    @{ Row = 134; A = 14386751; B = "2025-08-08"; C = "Cristina Bucsa";   D = "Yue Yuan";         E = "Gana Yue Yuan";        F = 2;    G = "Acierto"; H = 1 },
    @{ Row = 135; A = 14311731; B = "2025-08-09"; C = "Harold Mayot";     D = "Dino Prižmić";     E = "Gana Harold Mayot";    F = 3.25; G = $null;     H = $null },
    @{ Row = 136; A = 14311736; B = "2025-08-09"; C = "Kamil Majchrzak";  D = "Ugo Blanchet";     E = "Gana Ugo Blanchet";    F = 3;    G = $null;     H = $null }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A" + $row).Value = $r.A
    $ws.Range("B" + $row).NumberFormat = "@"
    $ws.Range("B" + $row).Value = $r.B
    $ws.Range("C" + $row).Value = $r.C
    $ws.Range("D" + $row).Value = $r.D
    $ws.Range("E" + $row).Value = $r.E
    $ws.Range("F" + $row).Value = $r.F
    if ($r.G -ne $null) {
        $ws.Range("G" + $row).Value = $r.G
    }
    if ($r.H -ne $null) {
        $ws.Range("H" + $row).Value = $r.H
    }
}
